$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.005.74"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.887.07"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.23"
$ws.Range("E5").Value = "  -2.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4599"
$ws.Range("E7").Value = "  -2.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4088"
$ws.Range("E8").Value = "  +2.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.30"
$ws.Range("E9").Value = "  -0.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07985"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9909"
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.71"
$ws.Range("E12").Value = "  -1.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.900.54"
$ws.Range("E13").Value = "  +1.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.901"
$ws.Range("E14").Value = "  -2.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.068"
$ws.Range("E15").Value = "  -2.94%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9995"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.85"
$ws.Range("E17").Value = "  -1.60%  "
$ws.Range("E18").Value = "  -1.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06556"
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9990"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.059.79"
$ws.Range("E22").Value = "  +1.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.402"
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.211"
$ws.Range("E25").Value = "  -1.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.106.80"
$ws.Range("E26").Value = "  +0.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.20"
$ws.Range("E27").Value = "  -2.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.62"
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.101"
$ws.Range("E29").Value = "  -1.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.406"
$ws.Range("E30").Value = "  -1.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.79"
$ws.Range("E31").Value = "  -1.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9759"
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09329"
$ws.Range("E33").Value = "  -2.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.597"
$ws.Range("E34").Value = "  -2.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.407"
$ws.Range("E35").Value = "  +1.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.273"
$ws.Range("E36").Value = "  -1.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06037"
$ws.Range("E37").Value = "  -2.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02223"
$ws.Range("E38").Value = "  -1.45%  "
$ws.Range("E39").Value = "  -2.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.179"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9986"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5772"
$ws.Range("E42").Value = "  -3.07%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1820"
$ws.Range("E43").Value = "  -3.41%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.11"
$ws.Range("E44").Value = "  -1.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.244"
$ws.Range("E45").Value = "  -1.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.263"
$ws.Range("E46").Value = "  +9.12%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.05"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5458"
$ws.Range("E48").Value = "  -2.07%  "
$ws.Range("E49").Value = "  -3.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07018"
$ws.Range("E50").Value = "  -5.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.45"
$ws.Range("E51").Value = "  +13.75%  "
